$d = $word.ActiveDocument

# 1. Rechtsgrundlage paragraph - rewrite sentence
$d.Content.Find.Execute(
    "Die Anordnung des Landratsamts könnte auf §1 Abs. 1 in Verbindung mit § 7 Abs.1 Satz 1 DSchG gestützt werden.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Rechtsgrundlage für die Anordnung der Reparatur des Fachwerkhauses mit Biberschwanz-Dachziegeln könnte §1 Abs. 1 in Verbindung mit § 7 Abs.1 Satz 1 DSchG sein.",
    2)

# 2. Der Pflichtige paragraph - insert new sentences about F.K before the G.K sentence
$d.Content.Find.Execute(
    "Als Pflichtige kommen sowohl F.K als auch G.K in Betracht. F.K könnte pflichtig sein sinngemäß § 7 Abs. 1 Satz 1 DSchG und § 7 PolG, da er Eigentümer einer Sache ist, von deren Zustand eine Gefahr ausgeht. Der G.K ist ebenfalls Eigentümer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Als Pflichtige kommen sowohl F.K als auch G.K in Betracht. F.K könnte pflichtig sein sinngemäß § 7 Abs. 1 Satz 1 DSchG und § 7 PolG, da er Eigentümer einer Sache ist, von deren Zustand eine Gefahr ausgeht. Der F.K ist Eigentümer des Fachwerkhauses, von dessen Dach eine Gefährdung für das Denkmal ausgeht. Also ist er Pflichtiger. Der G.K ist ebenfalls Eigentümer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.",
    2)

# 3. Sachliche Zuständigkeit paragraph - remove "§ " before "46 Abs. 2 LBO"
$d.Content.Find.Execute(
    "Nach §§ 7 Abs. 4, 3 Abs. 3, Abs. 1 Nr. 3 DSchG und § 46 Abs. 2 LBO und § 15 LVG ist das Landratsamt sachlich zuständig.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nach §§ 7 Abs. 4, 3 Abs. 3, Abs. 1 Nr. 3 DSchG und 46 Abs. 2 LBO und § 15 LVG ist das Landratsamt sachlich zuständig.",
    2)

# 4. Ausgeschlossene Personen/Befangenheit paragraph - rewrite sentence
$d.Content.Find.Execute(
    "Es liegen keine Hinweise auf eine Befangenheit vor.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Es gibt keine Hinweise auf ausgeschlossene Personen oder Befangenheit.",
    2)

# 5. Beteiligung anderer Behörden heading - drop trailing period
$d.Content.Find.Execute(
    "Beteiligung anderer Behörden.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Beteiligung anderer Behörden",
    2)
